$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 90
$ws.Cells.Item(90, 1).Value = 89
$ws.Cells.Item(90, 2).Value = 1
$ws.Cells.Item(90, 3).Value = "2024-06-16 19:11:54"
$ws.Cells.Item(90, 4).Value = 200
$ws.Cells.Item(90, 5).Value = 10

# Row 91
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = 2
$ws.Cells.Item(91, 3).Value = "2024-06-16 19:11:54"
$ws.Cells.Item(91, 4).Value = 200
$ws.Cells.Item(91, 5).Value = 2
